$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move wrapped-around "Doing" column content back from F/E into B/C,
#     matching the pre-merge layout (reverting the merge commit). ---

# Row 17: Passenger / Enter Passenger details -> column C
$ws.Range("F17").Copy($ws.Range("C17"))
$ws.Range("F17").Clear()

# Row 18: Payment (User) -> column B
$ws.Range("F18").Copy($ws.Range("B18"))
$ws.Range("F18").Clear()

# Row 22: Forgot Password (User) -> column C
$ws.Range("F22").Copy($ws.Range("C22"))
$ws.Range("F22").Clear()

# Row 23: Password Validation -> column C (keep source style s=4)
$ws.Range("F23").Copy($ws.Range("C23"))
$ws.Range("F23").Clear()

# Row 24: Admin View -> column C (keep source style s=4)
$ws.Range("E24").Copy($ws.Range("C24"))
$ws.Range("E24").Clear()

# Row 25: User Profile -> column B
$ws.Range("E25").Copy($ws.Range("B25"))
$ws.Range("E25").Clear()

# Rows 19-21, 26-28, 30-31 keep their existing A/B/C content; just make sure
# any stray D:F cells are vacated so the row's used-range (and therefore the
# serialized "spans") shrinks back down to column C.
$ws.Range("D19:F21").Clear()
$ws.Range("D26:F28").Clear()
$ws.Range("D30:F31").Clear()

# --- Row-height adjustments that came along with the narrower layout ---
$ws.Rows("18").RowHeight = 57.6
$ws.Rows("22").RowHeight = 28.8
$ws.Rows("24").RowHeight = 23.4
$ws.Rows("25").RowHeight = 18.6

# --- Restore the active selection/cursor position ---
$ws.Range("G14").Select()
